$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 1
$ws.Range("C62").Value = "2024-06-16 05:13:06"
$ws.Range("D62").Value = 200
$ws.Range("E62").Value = 4

$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 2
$ws.Range("C63").Value = "2024-06-16 05:13:06"
$ws.Range("D63").Value = 200
$ws.Range("E63").Value = 0
